# Update the HHS grant history table text to reflect FY 2012-2016
# instead of FY 2011-2016 (the data table itself already covers
# 2012-2016; only the descriptive text cells lagged behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HHSgrantHistTable")

$ws.Range("A3").Value = "This table shows the grant awards and award dollars HHS made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the HHS page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars HHS made for FY 2012-2016."
